$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1405
$ws.Range("F4").Value = 13552
$ws.Range("F5").Value = 789
$ws.Range("F8").Value = 65
$ws.Range("F10").Value = 1935
$ws.Range("F13").Value = 24686
$ws.Range("F14").Value = 551
$ws.Range("F15").Value = 230
$ws.Range("F16").Value = 537
$ws.Range("F17").Value = 146
$ws.Range("F19").Value = 230
$ws.Range("F20").Value = 331
$ws.Range("F21").Value = 177
$ws.Range("F23").Value = 38
$ws.Range("F25").Value = 299
$ws.Range("F26").Value = 28
$ws.Range("F27").Value = 1385
$ws.Range("F28").Value = 94
$ws.Range("F29").Value = 387
$ws.Range("F30").Value = 84

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 204
$ws.Range("F6").Value = 32
$ws.Range("F8").Value = 96
$ws.Range("F9").Value = 96
$ws.Range("F15").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 912
$ws.Range("F3").Value = 4657
$ws.Range("F4").Value = 134

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 912
$ws.Range("F4").Value = 1405
$ws.Range("F5").Value = 13552
$ws.Range("F6").Value = 789
$ws.Range("F7").Value = 4657
$ws.Range("F10").Value = 65
$ws.Range("F11").Value = 1935
$ws.Range("F13").Value = 134
$ws.Range("F14").Value = 24686
$ws.Range("F15").Value = 551
$ws.Range("F17").Value = 230
$ws.Range("F18").Value = 204
$ws.Range("F19").Value = 204
$ws.Range("F20").Value = 537
$ws.Range("F23").Value = 146
$ws.Range("F24").Value = 32
$ws.Range("F26").Value = 96
$ws.Range("F30").Value = 230
$ws.Range("F31").Value = 331
$ws.Range("F32").Value = 177
$ws.Range("F34").Value = 38
$ws.Range("F39").Value = 299
$ws.Range("F40").Value = 28
$ws.Range("F41").Value = 12
$ws.Range("F42").Value = 1385
$ws.Range("F43").Value = 94
$ws.Range("F45").Value = 387
$ws.Range("F46").Value = 84

